$d = $word.ActiveDocument

# Remove the hyperlink wrappers around "R will read the source files in
# alphabetical order" and "has a little-known feature" while keeping the
# run text and its formatting (Hyperlink.Delete removes the hyperlink
# field but leaves the display text/run in place).
for ($i = $d.Hyperlinks.Count; $i -ge 1; $i--) {
    $d.Hyperlinks.Item($i).Delete()
}
